# Delete the first data row (row 2) of the active sheet, shifting all
# subsequent rows up by one. This matches the author's commit of
# "Plotting GRS directly from work_vhe instead of copy on voids and GRS
# plot" which dropped the first record (A=88, 3FHL J0648.7+1517) from the
# table, shrinking the used range from A1:G122 down to A1:G121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
